$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Product " header had a trailing space -- trim it.
$ws.Range("D1").Value = "Product"

# Touch F3 (left-aligned, like the other input cells) so it becomes part of
# the used range as an empty, formatted input cell.
$xlLeft = -4131
$ws.Range("F3").HorizontalAlignment = $xlLeft

# Leave the cursor where the user left off after filling the form.
$ws.Range("E8").Select()
